$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# Row 2
$ws.Range("D2").Value = 117018
$ws.Range("E2").Value = 5799
$ws.Range("F2").Value = 5799
$ws.Range("G2").Value = 2430
$ws.Range("H2").Value = 1379
$ws.Range("I2").Value = 920
$ws.Range("J2").Value = 459
$ws.Range("K2").Value = 133825
$ws.Range("L2").Value = 84327
$ws.Range("M2").Value = 49499
$ws.Range("N2").Value = 29849
$ws.Range("O2").Value = 19650
$ws.Range("P2").Value = 723
$ws.Range("Q2").Value = 6886
$ws.Range("R2").Value = -5372
$ws.Range("S2").Value = -1414
$ws.Range("T2").Value = 4743
$ws.Range("U2").Value = 2143
$ws.Range("V2").Value = 59473
$ws.Range("W2").Value = 4.96
$ws.Range("X2").Value = 1.18
$ws.Range("Y2").Value = 3.13
$ws.Range("Z2").Value = 1.04
$ws.Range("AA2").Value = 170.36
$ws.Range("AB2").Value = 4288.09
$ws.Range("AC2").Value = 6369
$ws.Range("AD2").Value = 48.12
$ws.Range("AE2").Value = 211670
$ws.Range("AF2").Value = 1.45
$ws.Range("AG2").Value = 2000
$ws.Range("AH2").Value = 0.65
$ws.Range("AI2").Value = 30.74
$ws.Range("AJ2").Value = 13138298

# Row 3
$ws.Range("D3").Value = 129245
$ws.Range("E3").Value = 7514
$ws.Range("F3").Value = 7514
$ws.Range("G3").Value = 3648
$ws.Range("H3").Value = 2537
$ws.Range("I3").Value = 1893
$ws.Range("J3").Value = 644
$ws.Range("K3").Value = 137513
$ws.Range("L3").Value = 83660
$ws.Range("M3").Value = 53852
$ws.Range("N3").Value = 31431
$ws.Range("O3").Value = 22421
$ws.Range("P3").Value = 724
$ws.Range("Q3").Value = 11546
$ws.Range("R3").Value = -6930
$ws.Range("S3").Value = -3645
$ws.Range("T3").Value = 6735
$ws.Range("U3").Value = 4811
$ws.Range("V3").Value = 57462
$ws.Range("W3").Value = 5.81
$ws.Range("X3").Value = 1.96
$ws.Range("Y3").Value = 6.18
$ws.Range("Z3").Value = 1.87
$ws.Range("AA3").Value = 155.35
$ws.Range("AB3").Value = 4496.86
$ws.Range("AC3").Value = 13073
$ws.Range("AD3").Value = 28.88
$ws.Range("AE3").Value = 222486
$ws.Range("AF3").Value = 1.7
$ws.Range("AG3").Value = 2500
$ws.Range("AH3").Value = 0.66
$ws.Range("AI3").Value = 18.7
$ws.Range("AJ3").Value = 13158198

# Row 4
$ws.Range("D4").Value = 145633
$ws.Range("E4").Value = 8436
$ws.Range("F4").Value = 8436
$ws.Range("G4").Value = 5284
$ws.Range("H4").Value = 3535
$ws.Range("I4").Value = 2759
$ws.Range("J4").Value = 776
$ws.Range("K4").Value = 157662
$ws.Range("L4").Value = 97058
$ws.Range("M4").Value = 60605
$ws.Range("N4").Value = 34094
$ws.Range("O4").Value = 26511
$ws.Range("P4").Value = 724
$ws.Range("Q4").Value = 8001
$ws.Range("R4").Value = -14089
$ws.Range("S4").Value = 7259
$ws.Range("T4").Value = 8454
$ws.Range("U4").Value = -453
$ws.Range("V4").Value = 66622
$ws.Range("W4").Value = 5.79
$ws.Range("X4").Value = 2.43
$ws.Range("Y4").Value = 8.42
$ws.Range("Z4").Value = 2.4
$ws.Range("AA4").Value = 160.15
$ws.Range("AB4").Value = 4818.68
$ws.Range("AC4").Value = 19044
$ws.Range("AD4").Value = 18.77
$ws.Range("AE4").Value = 241270
$ws.Range("AF4").Value = 1.48
$ws.Range("AG4").Value = 2500
$ws.Range("AH4").Value = 0.7
$ws.Range("AI4").Value = 12.83
$ws.Range("AJ4").Value = 13161898

# Row 5
$ws.Range("D5").Value = 164772
$ws.Range("E5").Value = 7766
$ws.Range("F5").Value = 7766
$ws.Range("G5").Value = 5843
$ws.Range("H5").Value = 4128
$ws.Range("I5").Value = 3702
$ws.Range("J5").Value = 426
$ws.Range("K5").Value = 168681
$ws.Range("L5").Value = 107243
$ws.Range("M5").Value = 61439
$ws.Range("N5").Value = 33749
$ws.Range("O5").Value = 27689
$ws.Range("P5").Value = 725
$ws.Range("Q5").Value = 11808
$ws.Range("R5").Value = -14639
$ws.Range("S5").Value = 2654
$ws.Range("T5").Value = 14504
$ws.Range("U5").Value = -2695
$ws.Range("V5").Value = 71344
$ws.Range("W5").Value = 4.71
$ws.Range("X5").Value = 2.5
$ws.Range("Y5").Value = 10.91
$ws.Range("Z5").Value = 2.53
$ws.Range("AA5").Value = 174.55
$ws.Range("AB5").Value = 5152.2
$ws.Range("AC5").Value = 25536
$ws.Range("AD5").Value = 14.33
$ws.Range("AE5").Value = 238582
$ws.Range("AF5").Value = 1.53
$ws.Range("AG5").Value = 3000
$ws.Range("AH5").Value = 0.82
$ws.Range("AI5").Value = 11.48
$ws.Range("AJ5").Value = 13176748

# Row 6
$ws.Range("D6").Value = 186701
$ws.Range("E6").Value = 8327
$ws.Range("F6").Value = 8327
$ws.Range("G6").Value = 12942
$ws.Range("H6").Value = 9254
$ws.Range("I6").Value = 8752
$ws.Range("K6").Value = 194970
$ws.Range("L6").Value = 121888
$ws.Range("M6").Value = 73082
$ws.Range("N6").Value = 47828
$ws.Range("P6").Value = 819
$ws.Range("Q6").Value = 4714
$ws.Range("R6").Value = -10346
$ws.Range("S6").Value = 5210
$ws.Range("T6").Value = 16754
$ws.Range("U6").Value = -12040
$ws.Range("V6").Value = 79353
$ws.Range("W6").Value = 4.46
$ws.Range("X6").Value = 4.96
$ws.Range("Y6").Value = 21.46
$ws.Range("Z6").Value = 5.09
$ws.Range("AA6").Value = 166.78
$ws.Range("AB6").Value = 6233.74
$ws.Range("AC6").Value = 54173
$ws.Range("AD6").Value = 6.1
$ws.Range("AE6").Value = 298494
$ws.Range("AF6").Value = 1.11
$ws.Range("AG6").Value = 3500
$ws.Range("AH6").Value = 1.06
$ws.Range("AI6").Value = 6.42
$ws.Range("AJ6").Value = 15054186

# Row 7
$ws.Range("D7").Value = 223437
$ws.Range("E7").Value = 8202
$ws.Range("G7").Value = 3084
$ws.Range("H7").Value = 1920
$ws.Range("I7").Value = 1393
$ws.Range("K7").Value = 229456
$ws.Range("L7").Value = 147894
$ws.Range("M7").Value = 81562
$ws.Range("N7").Value = 49068
$ws.Range("P7").Value = 820
$ws.Range("Q7").Value = 14894
$ws.Range("R7").Value = -23540
$ws.Range("S7").Value = 8477
$ws.Range("T7").Value = 14524
$ws.Range("U7").Value = -2179
$ws.Range("W7").Value = 3.67
$ws.Range("X7").Value = 0.86
$ws.Range("Y7").Value = 2.87
$ws.Range("Z7").Value = 0.91
$ws.Range("AA7").Value = 181.33
$ws.Range("AC7").Value = 8501
$ws.Range("AD7").Value = 27.94
$ws.Range("AE7").Value = 306234
$ws.Range("AF7").Value = 0.78
$ws.Range("AG7").Value = 3506
$ws.Range("AH7").Value = 1.48
$ws.Range("AI7").Value = 37.9

# Row 8
$ws.Range("D8").Value = 244910
$ws.Range("E8").Value = 9816
$ws.Range("G8").Value = 5414
$ws.Range("H8").Value = 3806
$ws.Range("I8").Value = 3097
$ws.Range("K8").Value = 236124
$ws.Range("L8").Value = 151250
$ws.Range("M8").Value = 84874
$ws.Range("N8").Value = 51549
$ws.Range("P8").Value = 820
$ws.Range("Q8").Value = 14399
$ws.Range("R8").Value = -10499
$ws.Range("S8").Value = -3502
$ws.Range("T8").Value = 10540
$ws.Range("U8").Value = 4312
$ws.Range("W8").Value = 4.01
$ws.Range("X8").Value = 1.55
$ws.Range("Y8").Value = 6.16
$ws.Range("Z8").Value = 1.64
$ws.Range("AA8").Value = 178.21
$ws.Range("AC8").Value = 18906
$ws.Range("AD8").Value = 12.56
$ws.Range("AE8").Value = 321717
$ws.Range("AF8").Value = 0.74
$ws.Range("AG8").Value = 3562
$ws.Range("AH8").Value = 1.5
$ws.Range("AI8").Value = 17.32

# Row 9
$ws.Range("D9").Value = 263991
$ws.Range("E9").Value = 11078
$ws.Range("G9").Value = 6771
$ws.Range("H9").Value = 4709
$ws.Range("I9").Value = 3863
$ws.Range("K9").Value = 243767
$ws.Range("L9").Value = 154682
$ws.Range("M9").Value = 89084
$ws.Range("N9").Value = 54748
$ws.Range("P9").Value = 820
$ws.Range("Q9").Value = 14364
$ws.Range("R9").Value = -10687
$ws.Range("S9").Value = -2956
$ws.Range("T9").Value = 10339
$ws.Range("U9").Value = 4323
$ws.Range("W9").Value = 4.2
$ws.Range("X9").Value = 1.78
$ws.Range("Y9").Value = 7.27
$ws.Range("Z9").Value = 1.96
$ws.Range("AA9").Value = 173.64
$ws.Range("AC9").Value = 23580
$ws.Range("AD9").Value = 10.07
$ws.Range("AE9").Value = 341680
$ws.Range("AF9").Value = 0.7
$ws.Range("AG9").Value = 3706
$ws.Range("AH9").Value = 1.56
$ws.Range("AI9").Value = 14.44
